$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 (Thursday) rework ---
# Before: I5:N5 merged, holds "CS161 | Problem Solving through Programming | Dr. Sunil C K | C003" (style 4)
#         O5:U5 each hold "LUNCH BREAK" (style 2), not merged
# After:  I5:M5 cleared/plain, N5:Q5 merged "LUNCH BREAK" (style 2), R5:W5 merged new CS161 (style 4, room C002)

# Unmerge the old CS161 block at I5:N5 first
$ws.Range("I5:N5").UnMerge()

# Clear the old CS161 text/style from I5:N5 and the old individual LUNCH BREAK cells O5:U5
$ws.Range("I5:U5").ClearContents()
$ws.Range("I5:U5").Style = "Normal"

# New lunch break merged block N5:Q5
$ws.Range("N5:Q5").Merge()
$ws.Range("N5").Value = "LUNCH BREAK"
$ws.Range("N5:Q5").Style = $ws.Range("N4:Q4").Style
$ws.Range("N5:Q5").HorizontalAlignment = -4108
$ws.Range("N5:Q5").VerticalAlignment = -4108
$ws.Range("N5:Q5").WrapText = $true

# New CS161 merged block R5:W5 (room changed to C002)
$ws.Range("R5:W5").Merge()
$ws.Range("R5").Value = "CS161 | Problem Solving through Programming | Dr. Sunil C K | C002"
$ws.Range("R5:W5").Style = $ws.Range("R4:W4").Style
$ws.Range("R5:W5").HorizontalAlignment = -4108
$ws.Range("R5:W5").VerticalAlignment = -4108
$ws.Range("R5:W5").WrapText = $true
